$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 85.818184
$ws.Range("I9").Value = 85.818184
$ws.Range("K9").Value = 85.818184
$ws.Range("M9").Value = 83.181816
$ws.Range("H15").Value = 1071.64
$ws.Range("I15").Value = 1071.64
$ws.Range("K15").Value = 3214.92
$ws.Range("M15").Value = -3045.92
$ws.Range("H28").Value = 1079
$ws.Range("I28").Value = 895.8889
$ws.Range("K28").Value = 895.8889
$ws.Range("M28").Value = -410.8889
$ws.Range("H86").Value = 2872.625
$ws.Range("I86").Value = 2376.2
$ws.Range("J86").Value = 3700
$ws.Range("K86").Value = 2376.2
$ws.Range("L86").Value = 3700
$ws.Range("M86").Value = -1253.2
$ws.Range("N86").Value = -5946
$ws.Range("H88").Value = 3857.2856
$ws.Range("J88").Value = 6000.5
$ws.Range("L88").Value = 6000.5
$ws.Range("N88").Value = -6812.5
$ws.Range("H89").Value = 2872.625
$ws.Range("I89").Value = 2376.2
$ws.Range("J89").Value = 3700
$ws.Range("K89").Value = 11881
$ws.Range("L89").Value = 18500
$ws.Range("M89").Value = -6265
$ws.Range("N89").Value = -29732
$ws.Range("H91").Value = 3857.2856
$ws.Range("J91").Value = 6000.5
$ws.Range("L91").Value = 6000.5
$ws.Range("N91").Value = -8808.5
$ws.Range("H106").Value = 33761.832
$ws.Range("I106").Value = 36822
$ws.Range("J106").Value = 100
$ws.Range("K106").Value = 36822
$ws.Range("L106").Value = 100
$ws.Range("M106").Value = -36191
$ws.Range("N106").Value = -1362
$ws.Range("H111").Value = 2983.9333
$ws.Range("I111").Value = 1877.8889
$ws.Range("K111").Value = 5633.6667
$ws.Range("M111").Value = -2566.6667
$ws.Range("H116").Value = 8100
$ws.Range("J116").Value = 8725
$ws.Range("L116").Value = 8725
$ws.Range("N116").Value = -15609
$ws.Range("H132").Value = 2066.9375
$ws.Range("I132").Value = 2290.7856
$ws.Range("K132").Value = 6872.3568
$ws.Range("M132").Value = -4342.3568
$ws.Range("H137").Value = 2443.6365
$ws.Range("I137").Value = 1827.8462
$ws.Range("K137").Value = 5483.5386
$ws.Range("M137").Value = -2933.5386

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3250
$ws.Range("J2").Value = 5000
$ws.Range("L2").Value = 5000
$ws.Range("N2").Value = -5226
$ws.Range("H32").Value = 4646.4
$ws.Range("I32").Value = 2864.9524
$ws.Range("J32").Value = 13999
$ws.Range("K32").Value = 2864.9524
$ws.Range("L32").Value = 13999
$ws.Range("M32").Value = -2577.9524
$ws.Range("N32").Value = -14573
$ws.Range("H43").Value = 99999.5
$ws.Range("J43").Value = 99999.5
$ws.Range("L43").Value = 99999.5
$ws.Range("N43").Value = -100625.5
$ws.Range("H45").Value = 1860.4
$ws.Range("I45").Value = 1825.5
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1825.5
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1448.5
$ws.Range("N45").Value = -2754
$ws.Range("H74").Value = 1298.4
$ws.Range("I74").Value = 878.3125
$ws.Range("J74").Value = 5779.3335
$ws.Range("K74").Value = 878.3125
$ws.Range("L74").Value = 5779.3335
$ws.Range("M74").Value = -4.3125
$ws.Range("N74").Value = -7527.3335
$ws.Range("H77").Value = 1298.4
$ws.Range("I77").Value = 878.3125
$ws.Range("J77").Value = 5779.3335
$ws.Range("K77").Value = 4391.5625
$ws.Range("L77").Value = 28896.6675
$ws.Range("M77").Value = -23.5625
$ws.Range("N77").Value = -37632.6675
$ws.Range("H88").Value = 1176.5454
$ws.Range("I88").Value = 761.25
$ws.Range("J88").Value = 1413.8572
$ws.Range("K88").Value = 761.25
$ws.Range("L88").Value = 1413.8572
$ws.Range("M88").Value = -355.25
$ws.Range("N88").Value = -2225.8572
$ws.Range("H91").Value = 1176.5454
$ws.Range("I91").Value = 761.25
$ws.Range("J91").Value = 1413.8572
$ws.Range("K91").Value = 761.25
$ws.Range("L91").Value = 1413.8572
$ws.Range("M91").Value = 642.75
$ws.Range("N91").Value = -4221.8572
$ws.Range("H116").Value = 3250
$ws.Range("J116").Value = 5000
$ws.Range("L116").Value = 5000
$ws.Range("N116").Value = -9588
$ws.Range("H132").Value = 1576.4348
$ws.Range("I132").Value = 1475.4
$ws.Range("K132").Value = 4426.200000000001
$ws.Range("M132").Value = -1896.200000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3250
$ws.Range("J3").Value = 5000
$ws.Range("L3").Value = 5000
$ws.Range("N3").Value = -5228
$ws.Range("H134").Value = 2552.879
$ws.Range("I134").Value = 2403.0908
$ws.Range("J134").Value = 2852.4546
$ws.Range("K134").Value = 7209.2724
$ws.Range("L134").Value = 8557.363799999999
$ws.Range("M134").Value = -4674.2724
$ws.Range("N134").Value = -13627.3638

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 44999.5
$ws.Range("J68").Value = 44999.5
$ws.Range("L68").Value = 44999.5
$ws.Range("N68").Value = -46497.5
$ws.Range("H71").Value = 44999.5
$ws.Range("J71").Value = 44999.5
$ws.Range("L71").Value = 134998.5
$ws.Range("N71").Value = -142486.5
$ws.Range("H116").Value = 100000
$ws.Range("J116").Value = 100000
$ws.Range("L116").Value = 100000
$ws.Range("N116").Value = -109178
$ws.Range("H132").Value = 2302.682
$ws.Range("I132").Value = 1813.4474
$ws.Range("J132").Value = 5401.1665
$ws.Range("K132").Value = 5440.3422
$ws.Range("L132").Value = 16203.4995
$ws.Range("M132").Value = -2910.3422
$ws.Range("N132").Value = -21263.4995
$ws.Range("H134").Value = 2022.2449
$ws.Range("I134").Value = 1644.4166
$ws.Range("J134").Value = 3068.5386
$ws.Range("K134").Value = 4933.2498
$ws.Range("L134").Value = 9205.6158
$ws.Range("M134").Value = -2398.2498
$ws.Range("N134").Value = -14275.6158

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 237
$ws.Range("I50").Value = 237
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 711
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -230
$ws.Range("N50").ClearContents()
$ws.Range("H53").Value = 237
$ws.Range("I53").Value = 237
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 711
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -230
$ws.Range("N53").ClearContents()
$ws.Range("H68").Value = 300
$ws.Range("I68").Value = 300
$ws.Range("K68").Value = 900
$ws.Range("M68").Value = -89
$ws.Range("H71").Value = 300
$ws.Range("I71").Value = 300
$ws.Range("K71").Value = 2700
$ws.Range("M71").Value = 1356

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 12423.833
$ws.Range("J98").Value = 12423.833
$ws.Range("L98").Value = 12423.833
$ws.Range("N98").Value = -18413.833
$ws.Range("H107").Value = 127.5
$ws.Range("I107").Value = 120
$ws.Range("J107").Value = 135
$ws.Range("K107").Value = 120
$ws.Range("L107").Value = 135
$ws.Range("M107").Value = 1800
$ws.Range("N107").Value = -3975
$ws.Range("H114").Value = 89775
$ws.Range("J114").Value = 89775
$ws.Range("L114").Value = 89775
$ws.Range("N114").Value = -98453
$ws.Range("H118").Value = 24999.5
$ws.Range("J118").Value = 38999
$ws.Range("L118").Value = 38999
$ws.Range("N118").Value = -42313
$ws.Range("H132").Value = 2277.9062
$ws.Range("I132").Value = 1775.6
$ws.Range("J132").Value = 3115.0833
$ws.Range("K132").Value = 5326.799999999999
$ws.Range("L132").Value = 9345.249899999999
$ws.Range("M132").Value = -2796.799999999999
$ws.Range("N132").Value = -14405.2499

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H106").Value = 19399.8
$ws.Range("J106").Value = 19399.8
$ws.Range("L106").Value = 19399.8
$ws.Range("N106").Value = -21923.8
$ws.Range("H122").Value = 7826.222
$ws.Range("J122").Value = 5113.25
$ws.Range("L122").Value = 15339.75
$ws.Range("N122").Value = -20239.75
$ws.Range("H132").Value = 3326.879
$ws.Range("I132").Value = 2916.2083
$ws.Range("K132").Value = 8748.624899999999
$ws.Range("M132").Value = -6218.624899999999
$ws.Range("H136").Value = 3897.9
$ws.Range("I136").Value = 3533.2942
$ws.Range("K136").Value = 10599.8826
$ws.Range("M136").Value = -8049.882599999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 110000
$ws.Range("J99").Value = 70000
$ws.Range("L99").Value = 70000
$ws.Range("N99").Value = -75990
$ws.Range("H132").Value = 50812.527
$ws.Range("I132").Value = 63849.8
$ws.Range("K132").Value = 191549.4
$ws.Range("M132").Value = -189019.4
